# Update Name of Algo
# Apply updated KNN imputation results to the worksheet.
# The workbook contains a single data sheet ("Sheet1") laid out with
# headers A/B/C/D in row 1 and numeric data starting at row 2.
# This script writes the revised values for the specific cells that
# changed between the previous and current run of the algorithm.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -8.684999999999999
$ws.Range("D9").Value = -8.244000000000002
$ws.Range("A11").Value = -21.977
$ws.Range("C11").Value = -13.084
$ws.Range("A12").Value = -21.48800000000001
$ws.Range("D13").Value = -8.334
$ws.Range("D14").Value = -8.016999999999999
$ws.Range("A15").Value = -21.93
$ws.Range("D19").Value = -8.175999999999998
$ws.Range("D21").Value = -8.404
$ws.Range("D22").Value = -8.204000000000001
$ws.Range("C23").Value = -12.395
$ws.Range("D24").Value = -6.946000000000001
$ws.Range("D26").Value = -7.402000000000001
$ws.Range("A27").Value = -21.702
$ws.Range("A28").Value = -21.855
$ws.Range("C28").Value = -13.14
$ws.Range("A31").Value = -21.594
$ws.Range("A32").Value = -22.141
$ws.Range("C32").Value = -13.305
$ws.Range("C34").Value = -12.101
$ws.Range("A36").Value = -20.34
$ws.Range("C36").Value = -12.515
$ws.Range("C37").Value = -13.051
$ws.Range("A38").Value = -19.818
$ws.Range("D38").Value = -8.415000000000001
$ws.Range("D41").Value = -8.489000000000001
$ws.Range("C42").Value = -12.403
$ws.Range("A46").Value = -21.761
$ws.Range("C49").Value = -12.955
$ws.Range("D52").Value = -7.679
$ws.Range("A54").Value = -21.756
$ws.Range("C54").Value = -12.868
$ws.Range("A55").Value = -21.961
$ws.Range("A56").Value = -21.986
$ws.Range("D56").Value = -8.417999999999999
$ws.Range("A67").Value = -21.483
$ws.Range("A69").Value = -21.47
$ws.Range("D71").Value = -7.183000000000002
$ws.Range("A72").Value = -21.624
$ws.Range("D72").Value = -7.32
$ws.Range("A73").Value = -19.994
$ws.Range("C78").Value = -12.573
$ws.Range("D78").Value = -7.586
$ws.Range("C80").Value = -12.152
$ws.Range("A83").Value = -21.937
$ws.Range("D83").Value = -7.997999999999999
$ws.Range("D85").Value = -8.618
$ws.Range("A86").Value = -22.115
$ws.Range("D86").Value = -8.556000000000001
$ws.Range("D90").Value = -7.183
$ws.Range("A91").Value = -20.931
$ws.Range("A93").Value = -21.54
$ws.Range("D96").Value = -7.950999999999999
$ws.Range("C97").Value = -11.641
$ws.Range("A99").Value = -21.659
$ws.Range("C99").Value = -12.32
$ws.Range("C100").Value = -12.304
$ws.Range("C101").Value = -12.188
$ws.Range("D103").Value = -8.345000000000001
$ws.Range("A104").Value = -21.261
$ws.Range("A105").Value = -20.252
